# Generate Report for Handback
# - Overview sheet: flip the "Ready for handoff" status (shared by the
#   a6eefc6c row, both zh-cn and de-de columns) to "Handback transform failed".
# - zh-cn / de-de sheets: record the handback/handoff file-name mismatch in
#   the "Error Detail" column (K) for the a6eefc6c row (row 3).

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "Handback file name: bzd4if2i.az1 is different with handoff file name: a6eefc6c-4fad-4220-9f1f-e61abe03d3a9.b7776a49fb9a3b564ff416f4657fd24c76619435.zh-cn."

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "Handback file name: bzd4if2i.az1 is different with handoff file name: a6eefc6c-4fad-4220-9f1f-e61abe03d3a9.b7776a49fb9a3b564ff416f4657fd24c76619435.de-de."
